$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 33 (G33=5512)
$ws_ALC.Cells.Item(33, 8).Value = 168.54546
$ws_ALC.Cells.Item(33, 9).Value = 163.75
$ws_ALC.Cells.Item(33, 10).Value = 181.33333
$ws_ALC.Cells.Item(33, 11).Value = 163.75
$ws_ALC.Cells.Item(33, 12).Value = 181.33333
$ws_ALC.Cells.Item(33, 13).Value = 65.25
$ws_ALC.Cells.Item(33, 14).Value = -639.3333299999999

# ALC row 70 (G70=12604)
$ws_ALC.Cells.Item(70, 8).Value = 3909.0715
$ws_ALC.Cells.Item(70, 9).Value = 1547.2
$ws_ALC.Cells.Item(70, 10).Value = 5221.222
$ws_ALC.Cells.Item(70, 11).Value = 4641.6
$ws_ALC.Cells.Item(70, 12).Value = 15663.666
$ws_ALC.Cells.Item(70, 13).Value = -4371.6
$ws_ALC.Cells.Item(70, 14).Value = -16203.666

# ALC row 73 (G73=12604)
$ws_ALC.Cells.Item(73, 8).Value = 3909.0715
$ws_ALC.Cells.Item(73, 9).Value = 1547.2
$ws_ALC.Cells.Item(73, 10).Value = 5221.222
$ws_ALC.Cells.Item(73, 11).Value = 4641.6
$ws_ALC.Cells.Item(73, 12).Value = 15663.666
$ws_ALC.Cells.Item(73, 13).Value = -3705.6
$ws_ALC.Cells.Item(73, 14).Value = -17535.666

# ALC row 74 (G74=5507)
$ws_ALC.Cells.Item(74, 8).Value = 2552.7144
$ws_ALC.Cells.Item(74, 9).Value = 2644.8333
$ws_ALC.Cells.Item(74, 11).Value = 2644.8333
$ws_ALC.Cells.Item(74, 13).Value = -1708.8333

# ALC row 77 (G77=5507)
$ws_ALC.Cells.Item(77, 8).Value = 2552.7144
$ws_ALC.Cells.Item(77, 9).Value = 2644.8333
$ws_ALC.Cells.Item(77, 11).Value = 13224.1665
$ws_ALC.Cells.Item(77, 13).Value = -8544.166499999999

# ALC row 86 (G86=12603)
$ws_ALC.Cells.Item(86, 8).Value = 5586.75
$ws_ALC.Cells.Item(86, 9).Value = 5956.2856
$ws_ALC.Cells.Item(86, 10).Value = 3000
$ws_ALC.Cells.Item(86, 11).Value = 5956.2856
$ws_ALC.Cells.Item(86, 12).Value = 3000
$ws_ALC.Cells.Item(86, 13).Value = -4833.2856
$ws_ALC.Cells.Item(86, 14).Value = -5246

# ALC row 89 (G89=12603)
$ws_ALC.Cells.Item(89, 8).Value = 5586.75
$ws_ALC.Cells.Item(89, 9).Value = 5956.2856
$ws_ALC.Cells.Item(89, 10).Value = 3000
$ws_ALC.Cells.Item(89, 11).Value = 29781.428
$ws_ALC.Cells.Item(89, 12).Value = 15000
$ws_ALC.Cells.Item(89, 13).Value = -24165.428
$ws_ALC.Cells.Item(89, 14).Value = -26232

# ALC row 106 (G106=19903)
$ws_ALC.Cells.Item(106, 8).Value = 3873.2222
$ws_ALC.Cells.Item(106, 9).Value = 3940.75
$ws_ALC.Cells.Item(106, 11).Value = 3940.75
$ws_ALC.Cells.Item(106, 13).Value = -3309.75

# ALC row 141 (G141=44161)
$ws_ALC.Cells.Item(141, 8).Value = 3820
$ws_ALC.Cells.Item(141, 9).Value = 2750
$ws_ALC.Cells.Item(141, 11).Value = 8250
$ws_ALC.Cells.Item(141, 13).Value = -3070

# ARM row 2 (G2=27713)
$ws_ARM.Cells.Item(2, 8).Value = 1037.2858
$ws_ARM.Cells.Item(2, 9).Value = 1037.2858
$ws_ARM.Cells.Item(2, 11).Value = 1037.2858
$ws_ARM.Cells.Item(2, 13).Value = -924.2858000000001

# ARM row 61 (G61=43999)
$ws_ARM.Cells.Item(61, 8).Value = 3853.04
$ws_ARM.Cells.Item(61, 9).Value = 2872.7144
$ws_ARM.Cells.Item(61, 11).Value = 2872.7144
$ws_ARM.Cells.Item(61, 13).Value = -2660.7144

# ARM row 97 (G97=19941)
$ws_ARM.Cells.Item(97, 8).Value = 485.25
$ws_ARM.Cells.Item(97, 9).Value = 485.25
$ws_ARM.Cells.Item(97, 11).Value = 485.25
$ws_ARM.Cells.Item(97, 13).Value = 10.75

# ARM row 102 (G102=19945)
$ws_ARM.Cells.Item(102, 8).Value = 4501.1816
$ws_ARM.Cells.Item(102, 9).Value = 2216.1428
$ws_ARM.Cells.Item(102, 10).Value = 8500
$ws_ARM.Cells.Item(102, 11).Value = 2216.1428
$ws_ARM.Cells.Item(102, 12).Value = 8500
$ws_ARM.Cells.Item(102, 13).Value = -594.1428000000001
$ws_ARM.Cells.Item(102, 14).Value = -11744

# ARM row 116 (G116=27713)
$ws_ARM.Cells.Item(116, 8).Value = 1037.2858
$ws_ARM.Cells.Item(116, 9).Value = 1037.2858
$ws_ARM.Cells.Item(116, 11).Value = 1037.2858
$ws_ARM.Cells.Item(116, 13).Value = 1256.7142

# ARM row 122 (G122=36168)
$ws_ARM.Cells.Item(122, 8).Value = 524.4545000000001
$ws_ARM.Cells.Item(122, 9).Value = 524.4545000000001
$ws_ARM.Cells.Item(122, 11).Value = 1573.3635
$ws_ARM.Cells.Item(122, 13).Value = 876.6364999999998

# ARM row 136 (G136=43999)
$ws_ARM.Cells.Item(136, 8).Value = 3853.04
$ws_ARM.Cells.Item(136, 9).Value = 2872.7144
$ws_ARM.Cells.Item(136, 11).Value = 8618.143199999999
$ws_ARM.Cells.Item(136, 13).Value = -6068.143199999999

# BSM row 3 (G3=27713)
$ws_BSM.Cells.Item(3, 8).Value = 1037.2858
$ws_BSM.Cells.Item(3, 9).Value = 1037.2858
$ws_BSM.Cells.Item(3, 11).Value = 1037.2858
$ws_BSM.Cells.Item(3, 13).Value = -923.2858000000001

# BSM row 20 (G20=14149)
$ws_BSM.Cells.Item(20, 8).Value = 2373.5625
$ws_BSM.Cells.Item(20, 9).Value = 1783.7142
$ws_BSM.Cells.Item(20, 10).Value = 2832.3333
$ws_BSM.Cells.Item(20, 11).Value = 1783.7142
$ws_BSM.Cells.Item(20, 12).Value = 2832.3333
$ws_BSM.Cells.Item(20, 13).Value = -1536.7142
$ws_BSM.Cells.Item(20, 14).Value = -3326.3333

# BSM row 22 (G22=5092)
$ws_BSM.Cells.Item(22, 8).Value = 1545.7142
$ws_BSM.Cells.Item(22, 9).Value = 1658.8462
$ws_BSM.Cells.Item(22, 10).Value = 75
$ws_BSM.Cells.Item(22, 11).Value = 1658.8462
$ws_BSM.Cells.Item(22, 12).Value = 75
$ws_BSM.Cells.Item(22, 13).Value = -1485.8462
$ws_BSM.Cells.Item(22, 14).Value = -421

# BSM row 80 (G80=13747)
$ws_BSM.Cells.Item(80, 8).Value = 1850
$ws_BSM.Cells.Item(80, 9).Value = 250
$ws_BSM.Cells.Item(80, 10).Value = 2650
$ws_BSM.Cells.Item(80, 11).Value = 250
$ws_BSM.Cells.Item(80, 12).Value = 2650
$ws_BSM.Cells.Item(80, 13).Value = 748
$ws_BSM.Cells.Item(80, 14).Value = -4646

# BSM row 83 (G83=13747)
$ws_BSM.Cells.Item(83, 8).Value = 1850
$ws_BSM.Cells.Item(83, 9).Value = 250
$ws_BSM.Cells.Item(83, 10).Value = 2650
$ws_BSM.Cells.Item(83, 11).Value = 1250
$ws_BSM.Cells.Item(83, 12).Value = 13250
$ws_BSM.Cells.Item(83, 13).Value = 3742
$ws_BSM.Cells.Item(83, 14).Value = -23234

# BSM row 86 (G86=12526)
$ws_BSM.Cells.Item(86, 8).Value = 3530.0833
$ws_BSM.Cells.Item(86, 9).Value = 1805.25
$ws_BSM.Cells.Item(86, 10).Value = 5254.9165
$ws_BSM.Cells.Item(86, 11).Value = 1805.25
$ws_BSM.Cells.Item(86, 12).Value = 5254.9165
$ws_BSM.Cells.Item(86, 13).Value = -682.25
$ws_BSM.Cells.Item(86, 14).Value = -7500.9165

# BSM row 89 (G89=12526)
$ws_BSM.Cells.Item(89, 8).Value = 3530.0833
$ws_BSM.Cells.Item(89, 9).Value = 1805.25
$ws_BSM.Cells.Item(89, 10).Value = 5254.9165
$ws_BSM.Cells.Item(89, 11).Value = 9026.25
$ws_BSM.Cells.Item(89, 12).Value = 26274.5825
$ws_BSM.Cells.Item(89, 13).Value = -3410.25
$ws_BSM.Cells.Item(89, 14).Value = -37506.5825

# BSM row 94 (G94=19939)
$ws_BSM.Cells.Item(94, 8).Value = 209.88889
$ws_BSM.Cells.Item(94, 9).Value = 209.88889
$ws_BSM.Cells.Item(94, 11).Value = 209.88889
$ws_BSM.Cells.Item(94, 13).Value = 241.11111

# BSM row 99 (G99=19943)
$ws_BSM.Cells.Item(99, 8).Value = 2093.0527
$ws_BSM.Cells.Item(99, 9).Value = 1561.1538
$ws_BSM.Cells.Item(99, 11).Value = 1561.1538
$ws_BSM.Cells.Item(99, 13).Value = -63.15380000000005

# BSM row 107 (G107=27706)
$ws_BSM.Cells.Item(107, 8).Value = 3506.4285
$ws_BSM.Cells.Item(107, 9).Value = 2424.1667
$ws_BSM.Cells.Item(107, 11).Value = 2424.1667
$ws_BSM.Cells.Item(107, 13).Value = -504.1667000000002

# CRP row 22 (G22=5367)
$ws_CRP.Cells.Item(22, 8).Value = 981.625
$ws_CRP.Cells.Item(22, 9).Value = 808.3333
$ws_CRP.Cells.Item(22, 11).Value = 808.3333
$ws_CRP.Cells.Item(22, 13).Value = -458.3333

# CRP row 31 (G31=44023)
$ws_CRP.Cells.Item(31, 8).Value = 4283.116
$ws_CRP.Cells.Item(31, 9).Value = 1973.5
$ws_CRP.Cells.Item(31, 11).Value = 1973.5
$ws_CRP.Cells.Item(31, 13).Value = -1678.5

# CRP row 34 (G34=44023)
$ws_CRP.Cells.Item(34, 8).Value = 4283.116
$ws_CRP.Cells.Item(34, 9).Value = 1973.5
$ws_CRP.Cells.Item(34, 11).Value = 1973.5
$ws_CRP.Cells.Item(34, 13).Value = -1771.5

# CRP row 122 (G122=36196)
$ws_CRP.Cells.Item(122, 8).Value = 1050.5
$ws_CRP.Cells.Item(122, 9).Value = 1050.5
$ws_CRP.Cells.Item(122, 11).Value = 3151.5
$ws_CRP.Cells.Item(122, 13).Value = -701.5

# CUL row 7 (G7=4728)
$ws_CUL.Cells.Item(7, 8).Value = 200
$ws_CUL.Cells.Item(7, 10).Value = 200
$ws_CUL.Cells.Item(7, 12).Value = 600
$ws_CUL.Cells.Item(7, 14).Value = -824

# CUL row 122 (G122=36078)
$ws_CUL.Cells.Item(122, 8).Value = 1049.75
$ws_CUL.Cells.Item(122, 9).Value = 747
$ws_CUL.Cells.Item(122, 10).Value = 1352.5
$ws_CUL.Cells.Item(122, 11).Value = 6723
$ws_CUL.Cells.Item(122, 12).Value = 12172.5
$ws_CUL.Cells.Item(122, 13).Value = -4273
$ws_CUL.Cells.Item(122, 14).Value = -17072.5

# CUL row 138 (G138=44105)
$ws_CUL.Cells.Item(138, 8).Value = 5726.75
$ws_CUL.Cells.Item(138, 9).Value = 1816
$ws_CUL.Cells.Item(138, 11).Value = 5448
$ws_CUL.Cells.Item(138, 13).Value = -308

# GSM row 102 (G102=36169)
$ws_GSM.Cells.Item(102, 8).Value = 3104.2144
$ws_GSM.Cells.Item(102, 9).Value = 2958.3845
$ws_GSM.Cells.Item(102, 11).Value = 2958.3845
$ws_GSM.Cells.Item(102, 13).Value = -1336.3845

# GSM row 126 (G126=36184)
$ws_GSM.Cells.Item(126, 8).Value = 1999.75
$ws_GSM.Cells.Item(126, 9).Value = 2333
$ws_GSM.Cells.Item(126, 11).Value = 6999
$ws_GSM.Cells.Item(126, 13).Value = -4529

# GSM row 132 (G132=44008)
$ws_GSM.Cells.Item(132, 8).Value = 2123.875
$ws_GSM.Cells.Item(132, 9).Value = 2213
$ws_GSM.Cells.Item(132, 11).Value = 6639
$ws_GSM.Cells.Item(132, 13).Value = -4109

# LTW row 16 (G16=5289)
$ws_LTW.Cells.Item(16, 8).Value = 1380
$ws_LTW.Cells.Item(16, 9).Value = 1276
$ws_LTW.Cells.Item(16, 10).Value = 1900
$ws_LTW.Cells.Item(16, 11).Value = 1276
$ws_LTW.Cells.Item(16, 12).Value = 1900
$ws_LTW.Cells.Item(16, 13).Value = -1106
$ws_LTW.Cells.Item(16, 14).Value = -2240

# LTW row 22 (G22=5277)
$ws_LTW.Cells.Item(22, 8).Value = 1484.5
$ws_LTW.Cells.Item(22, 9).Value = 1500
$ws_LTW.Cells.Item(22, 10).Value = 1479.3334
$ws_LTW.Cells.Item(22, 11).Value = 1500
$ws_LTW.Cells.Item(22, 12).Value = 1479.3334
$ws_LTW.Cells.Item(22, 13).Value = -1205
$ws_LTW.Cells.Item(22, 14).Value = -2069.3334

# LTW row 27 (G27=5277)
$ws_LTW.Cells.Item(27, 8).Value = 1484.5
$ws_LTW.Cells.Item(27, 9).Value = 1500
$ws_LTW.Cells.Item(27, 10).Value = 1479.3334
$ws_LTW.Cells.Item(27, 11).Value = 1500
$ws_LTW.Cells.Item(27, 12).Value = 1479.3334
$ws_LTW.Cells.Item(27, 13).Value = -1393
$ws_LTW.Cells.Item(27, 14).Value = -1693.3334

# LTW row 32 (G32=2250)
$ws_LTW.Cells.Item(32, 8).Value = 20013
$ws_LTW.Cells.Item(32, 9).Value = 20013
$ws_LTW.Cells.Item(32, 11).Value = 20013
$ws_LTW.Cells.Item(32, 13).Value = -19696

# LTW row 46 (G46=5282)
$ws_LTW.Cells.Item(46, 8).Value = 6259.875
$ws_LTW.Cells.Item(46, 10).Value = 6725.5713
$ws_LTW.Cells.Item(46, 12).Value = 6725.5713
$ws_LTW.Cells.Item(46, 14).Value = -7101.5713

# LTW row 68 (G68=12563)
$ws_LTW.Cells.Item(68, 8).Value = 6303.6665
$ws_LTW.Cells.Item(68, 10).Value = 8100
$ws_LTW.Cells.Item(68, 12).Value = 8100
$ws_LTW.Cells.Item(68, 14).Value = -9598

# LTW row 71 (G71=12563)
$ws_LTW.Cells.Item(71, 8).Value = 6303.6665
$ws_LTW.Cells.Item(71, 10).Value = 8100
$ws_LTW.Cells.Item(71, 12).Value = 40500
$ws_LTW.Cells.Item(71, 14).Value = -47988

# LTW row 93 (G93=19993)
$ws_LTW.Cells.Item(93, 8).Value = 1366.8
$ws_LTW.Cells.Item(93, 9).Value = 1038.8572
$ws_LTW.Cells.Item(93, 10).Value = 2132
$ws_LTW.Cells.Item(93, 11).Value = 1038.8572
$ws_LTW.Cells.Item(93, 12).Value = 2132
$ws_LTW.Cells.Item(93, 13).Value = 209.1428000000001
$ws_LTW.Cells.Item(93, 14).Value = -4628

# LTW row 100 (G100=19995)
$ws_LTW.Cells.Item(100, 8).Value = 5770.143
$ws_LTW.Cells.Item(100, 9).Value = 2540.4285
$ws_LTW.Cells.Item(100, 11).Value = 2540.4285
$ws_LTW.Cells.Item(100, 13).Value = -1999.4285

# LTW row 122 (G122=36247)
$ws_LTW.Cells.Item(122, 8).Value = 2505.7778
$ws_LTW.Cells.Item(122, 10).Value = 2776.7
$ws_LTW.Cells.Item(122, 12).Value = 8330.099999999999
$ws_LTW.Cells.Item(122, 14).Value = -13230.1

# LTW row 132 (G132=44058)
$ws_LTW.Cells.Item(132, 8).Value = 4314.8335
$ws_LTW.Cells.Item(132, 10).Value = 6502.5
$ws_LTW.Cells.Item(132, 12).Value = 19507.5
$ws_LTW.Cells.Item(132, 14).Value = -24567.5

# LTW row 138 (G138=42334)
$ws_LTW.Cells.Item(138, 8).Value = 80000
$ws_LTW.Cells.Item(138, 10).Value = 80000
$ws_LTW.Cells.Item(138, 12).Value = 80000
$ws_LTW.Cells.Item(138, 14).Value = -90280

# WVR row 100 (G100=19981)
$ws_WVR.Cells.Item(100, 8).Value = 555.5714
$ws_WVR.Cells.Item(100, 9).Value = 555.5714
$ws_WVR.Cells.Item(100, 11).Value = 1111.1428
$ws_WVR.Cells.Item(100, 13).Value = -570.1428000000001

# WVR row 122 (G122=36208)
$ws_WVR.Cells.Item(122, 8).Value = 4666.1665
$ws_WVR.Cells.Item(122, 9).Value = 2499.25
$ws_WVR.Cells.Item(122, 10).Value = 9000
$ws_WVR.Cells.Item(122, 11).Value = 7497.75
$ws_WVR.Cells.Item(122, 12).Value = 27000
$ws_WVR.Cells.Item(122, 13).Value = -5047.75
$ws_WVR.Cells.Item(122, 14).Value = -31900
